$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5604
$ws.Range("H3").Value = 8346
$ws.Range("I3").Value = 5850
$ws.Range("E4").Value = 1968
$ws.Range("H4").Value = 1674
$ws.Range("G5").Value = 788
$ws.Range("H5").Value = 803
$ws.Range("H6").Value = 7919
$ws.Range("I6").Value = 6495
$ws.Range("E7").Value = 25972
$ws.Range("G7").Value = 24668
$ws.Range("H7").Value = 25983
$ws.Range("I7").Value = 19809

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 56
$ws.Range("I7").Value = 230

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 60
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 209

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 196
$ws.Range("I3").Value = 211
$ws.Range("I6").Value = 178
$ws.Range("I7").Value = 634

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I6").Value = 93
$ws.Range("I7").Value = 366

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 288
$ws.Range("I7").Value = 784

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 169

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I6").Value = 130
$ws.Range("I7").Value = 460

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I5").Value = 63
$ws.Range("I6").Value = 146
$ws.Range("I7").Value = 629
$ws.Range("I8").Value = 1191
$ws.Range("I12").Value = 45
$ws.Range("I14").Value = 114
$ws.Range("I15").Value = 224
$ws.Range("I19").Value = 542
$ws.Range("I20").Value = 470
$ws.Range("I23").Value = 198
$ws.Range("I29").Value = 1245
$ws.Range("I31").Value = 195
$ws.Range("I36").Value = 259
$ws.Range("I37").Value = 634
$ws.Range("E42").Value = 1054
$ws.Range("I42").Value = 657
$ws.Range("I46").Value = 42
$ws.Range("I47").Value = 136
$ws.Range("I50").Value = 95
$ws.Range("I51").Value = 233
$ws.Range("H52").Value = 524
$ws.Range("I52").Value = 435
$ws.Range("I53").Value = 207
$ws.Range("I54").Value = 416
$ws.Range("I57").Value = 79
$ws.Range("I59").Value = 33
$ws.Range("G63").Value = 209
$ws.Range("H63").Value = 222
$ws.Range("I63").Value = 66
$ws.Range("I65").Value = 460
$ws.Range("I67").Value = 784
$ws.Range("I70").Value = 34
$ws.Range("I75").Value = 61
$ws.Range("I76").Value = 291
$ws.Range("I78").Value = 275
$ws.Range("I84").Value = 169
$ws.Range("I85").Value = 892
$ws.Range("I89").Value = 230
$ws.Range("I90").Value = 243
$ws.Range("I91").Value = 215
$ws.Range("I94").Value = 206
$ws.Range("I96").Value = 209
$ws.Range("I97").Value = 163
$ws.Range("I99").Value = 366
$ws.Range("E101").Value = 25972
$ws.Range("G101").Value = 24668
$ws.Range("H101").Value = 25983
$ws.Range("I101").Value = 19809

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 196
$ws.Range("I7").Value = 416

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 365
$ws.Range("I7").Value = 1245

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 165
$ws.Range("I7").Value = 542

$ws = $wb.Worksheets.Item("River North")
$ws.Range("H3").Value = 84
$ws.Range("H6").Value = 172
$ws.Range("I6").Value = 133
$ws.Range("I7").Value = 291

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I6").Value = 226
$ws.Range("I7").Value = 892

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 171
$ws.Range("E4").Value = 56
$ws.Range("I4").Value = 48
$ws.Range("E7").Value = 1054
$ws.Range("I7").Value = 657

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I6").Value = 102
$ws.Range("I7").Value = 275

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I3").Value = 79
$ws.Range("I7").Value = 215

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I6").Value = 152
$ws.Range("I7").Value = 470

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 77
$ws.Range("I6").Value = 78
$ws.Range("I7").Value = 259

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 118
$ws.Range("H4").Value = 28
$ws.Range("H7").Value = 524
$ws.Range("I7").Value = 435

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I3").Value = 34
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 206

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 136

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I2").Value = 69
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 224

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 26
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 33

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I3").Value = 10

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("I7").Value = 34

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I6").Value = 381
$ws.Range("I7").Value = 1191

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 243

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 60
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I4").Value = 8
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 44
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 207

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 209
$ws.Range("I3").Value = 195
$ws.Range("I7").Value = 629

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("I3").Value = 6
$ws.Range("I7").Value = 45
